$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header values for columns P (14) and Q (15) on row 1
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# For each data row (2-25): swap I<->K and M<->O, and add P=2, Q=2
for ($r = 2; $r -le 25; $r++) {
    $iVal = $ws.Cells.Item($r, 9).Value2   # column I
    $kVal = $ws.Cells.Item($r, 11).Value2  # column K
    $mVal = $ws.Cells.Item($r, 13).Value2  # column M
    $oVal = $ws.Cells.Item($r, 15).Value2  # column O

    $ws.Cells.Item($r, 9).Value = $kVal   # I = old K
    $ws.Cells.Item($r, 11).Value = $iVal  # K = old I
    $ws.Cells.Item($r, 13).Value = $oVal  # M = old O
    $ws.Cells.Item($r, 15).Value = $mVal  # O = old M

    $ws.Cells.Item($r, 16).Value = 2      # P
    $ws.Cells.Item($r, 17).Value = 2      # Q
}
